# Add "2022-Q4" sheet + update the "总计" (Total) summary sheet.
#
# Target layout after edit:
#   总计 (sheet 1, unchanged name/position)
#   2022-Q4 (NEW, inserted right after 总计)
#   2022-Q3 (previously sheet 2, now shifted to sheet 3)
#   2022-Q1 (previously sheet 3, now shifted to sheet 4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q4" sheet by duplicating "2022-Q3" (same column
#    layout/header/styles) and then overwriting its data rows.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")

$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# -- Row 2: fund 090019 / 大成景恒混合A --------------------------------
$q4Sheet.Range("B2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "090019"
$q4Sheet.Range("B2").ClearFormats()

$q4Sheet.Range("C2").NumberFormat = "@"
$q4Sheet.Range("C2").Value = "大成景恒混合A"
$q4Sheet.Range("C2").ClearFormats()

$q4Sheet.Range("D2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "1.18"
$q4Sheet.Range("D2").ClearFormats()

$q4Sheet.Range("E2").NumberFormat = "@"
$q4Sheet.Range("E2").Value = "93.72"
$q4Sheet.Range("E2").ClearFormats()

$q4Sheet.Range("F2").NumberFormat = "@"
$q4Sheet.Range("F2").Value = "1.48"
$q4Sheet.Range("F2").ClearFormats()

$q4Sheet.Range("G2").NumberFormat = "@"
$q4Sheet.Range("G2").Value = "0.0175"
$q4Sheet.Range("G2").ClearFormats()

$q4Sheet.Range("H2").Value = 10

# -- Row 3 (new row): fund 006038 / 大成景恒混合C -----------------------
# Clone A2's style (bold / bordered) onto A3 before filling it in.
$q4Sheet.Range("A2").Copy()
$q4Sheet.Range("A3").PasteSpecial(-4122)
$q4Sheet.Range("A3").Value = 1

$q4Sheet.Range("B3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "006038"
$q4Sheet.Range("B3").ClearFormats()

$q4Sheet.Range("C3").NumberFormat = "@"
$q4Sheet.Range("C3").Value = "大成景恒混合C"
$q4Sheet.Range("C3").ClearFormats()

$q4Sheet.Range("D3").NumberFormat = "@"
$q4Sheet.Range("D3").Value = "0.89"
$q4Sheet.Range("D3").ClearFormats()

$q4Sheet.Range("E3").NumberFormat = "@"
$q4Sheet.Range("E3").Value = "93.72"
$q4Sheet.Range("E3").ClearFormats()

$q4Sheet.Range("F3").NumberFormat = "@"
$q4Sheet.Range("F3").Value = "1.48"
$q4Sheet.Range("F3").ClearFormats()

$q4Sheet.Range("G3").NumberFormat = "@"
$q4Sheet.Range("G3").Value = "0.0132"
$q4Sheet.Range("G3").ClearFormats()

$q4Sheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a 2022-Q4 summary row right after
#    the header, shifting the old 2022-Q3 / 2022-Q1 rows down by one.
# ---------------------------------------------------------------------
# Shift old row 3 (2022-Q1) down to row 4 first (values, then formats).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4163)
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

# Shift old row 2 (2022-Q3) down to row 3 (values, then formats).
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4163)
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# Write the new 2022-Q4 summary into row 2.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

# Row indices (column A) need to read 0,1,2 top to bottom.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
